# Handback status refresh: swap the two tracked file UUIDs for a new pair
# and bump the generate/handback timestamps, matching a fresh CI run.
#
# Old UUID "4cc237eb-894e-46d4-a753-e8edc05a87a2" -> New "0b66ceda-aa33-44f6-ac71-a324e5d52ccc"
# Old UUID "f0496682-930e-4021-af04-f1d1ae7cd233" -> New "ffffde54a6f5-4b63-43a2-9d26-9a02e15157c5"
# Old xlf hash (both langs) "21bee339a012b50bf842e42ad59fa4c786adcc14" -> New "d45b6a77a8c23f294236e0a2bc4a44c8339e9b1e"
# Note: in the new run both rows collapse onto a single xlf name per language.

$wb = $excel.ActiveWorkbook

$oldUuid1 = "4cc237eb-894e-46d4-a753-e8edc05a87a2"
$oldUuid2 = "f0496682-930e-4021-af04-f1d1ae7cd233"
$newUuid1 = "0b66ceda-aa33-44f6-ac71-a324e5d52ccc"
$newUuid2 = "ffffde54a6f5-4b63-43a2-9d26-9a02e15157c5"

$newXlfZhCn = "$newUuid1.d45b6a77a8c23f294236e0a2bc4a44c8339e9b1e.zh-cn.xlf"
$newXlfDeDe = "$newUuid1.d45b6a77a8c23f294236e0a2bc4a44c8339e9b1e.de-de.xlf"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value2 = "$newUuid1.md"
$wsOverview.Range("B2").Value2 = "e2e\$newUuid1.md"
$wsOverview.Range("G2").Value2 = "2016-08-17 09:01:53"

$wsOverview.Range("A3").Value2 = "$newUuid2.md"
$wsOverview.Range("B3").Value2 = "e2e\$newUuid2.md"
$wsOverview.Range("G3").Value2 = "2016-08-17 09:01:53"

# Rebuild the two hyperlinks on this sheet (the host treats
# Hyperlinks.Item(n).Property = value as a no-op on the stored link, so the
# existing links are dropped and recreated with the new address/display).
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/be445d18fc0a4596efd3df3a9e4757636a5e8025/e2e/$newUuid1.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "e2e\$newUuid1.md"
)
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/be445d18fc0a4596efd3df3a9e4757636a5e8025/e2e/$newUuid2.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "e2e\$newUuid2.md"
)

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value2 = "$newUuid1.md"
$wsZhCn.Range("G2").Value2 = $newXlfZhCn
$wsZhCn.Range("H2").Value2 = "2016-08-17 09:01:48"
$wsZhCn.Range("I2").Value2 = "$newUuid1.md"
$wsZhCn.Range("J2").Value2 = $newXlfZhCn
$wsZhCn.Range("K2").Value2 = "2016-08-17 09:02:17"

$wsZhCn.Range("A3").Value2 = "$newUuid2.md"
$wsZhCn.Range("G3").Value2 = $newXlfZhCn
$wsZhCn.Range("H3").Value2 = "2016-08-17 09:01:48"
$wsZhCn.Range("I3").Value2 = "$newUuid2.md"
$wsZhCn.Range("J3").Value2 = $newXlfZhCn
$wsZhCn.Range("K3").Value2 = "2016-08-17 09:02:17"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/be445d18fc0a4596efd3df3a9e4757636a5e8025/e2e/$newUuid1.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "$newUuid1.md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/1abf6aec8d9d29949c12ab63e6740a714ca51e96/e2e/$newUuid1.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "$newUuid1.md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/be445d18fc0a4596efd3df3a9e4757636a5e8025/e2e/$newUuid2.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "$newUuid2.md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/1abf6aec8d9d29949c12ab63e6740a714ca51e96/e2e/$newUuid2.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "$newUuid2.md"
)

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value2 = "$newUuid1.md"
$wsDeDe.Range("G2").Value2 = $newXlfDeDe
$wsDeDe.Range("H2").Value2 = "2016-08-17 09:01:53"
$wsDeDe.Range("I2").Value2 = "$newUuid1.md"
$wsDeDe.Range("J2").Value2 = $newXlfDeDe
$wsDeDe.Range("K2").Value2 = "2016-08-17 09:02:25"

$wsDeDe.Range("A3").Value2 = "$newUuid2.md"
$wsDeDe.Range("G3").Value2 = $newXlfDeDe
$wsDeDe.Range("H3").Value2 = "2016-08-17 09:01:53"
$wsDeDe.Range("I3").Value2 = "$newUuid2.md"
$wsDeDe.Range("J3").Value2 = $newXlfDeDe
$wsDeDe.Range("K3").Value2 = "2016-08-17 09:02:25"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/be445d18fc0a4596efd3df3a9e4757636a5e8025/e2e/$newUuid1.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "$newUuid1.md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/343bb04a331c578b8d3de499f5e0f3e6bd2b358b/e2e/$newUuid1.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "$newUuid1.md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/be445d18fc0a4596efd3df3a9e4757636a5e8025/e2e/$newUuid2.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "$newUuid2.md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/343bb04a331c578b8d3de499f5e0f3e6bd2b358b/e2e/$newUuid2.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "$newUuid2.md"
)

Write-Host "Handback status refreshed."
